# Regenerate the "K" column (column G) values for the save_data sheet.
# The commit message indicates the underlying save-data generation logic
# changed ("use K instead of Strike#", recompute std/mean, etc.) and the
# K column values were recalculated and rewritten for each row. We apply
# the recalculated values directly to column G, matching the canonical
# OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 0
    8  = 4
    9  = 1
    10 = 3
    11 = 0
    12 = 0
    13 = 2
    14 = 0
    15 = 0
    16 = 2
    17 = 2
    18 = 3
    19 = 1
    20 = 0
    21 = 1
    22 = 1
    23 = 2
    25 = 2
    26 = 1
    27 = 1
    29 = 2
    31 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
